$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 33
$ws.Range("A33").Value = 7
$ws.Range('B33').Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range('C33').Value = 'Ñuble'
$ws.Range("D33").Value = (Get-Date -Year 2023 -Month 6 -Day 16 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E33").Value = 16
$ws.Range('F33').Value = 'Fruta'
$ws.Range("G33").Value = 100104
$ws.Range('H33').Value = 'Frutos de pepita'
$ws.Range("I33").Value = 100104003
$ws.Range('J33').Value = 'Membrillo'
$ws.Range('K33').Value = 'Champion'
$ws.Range('L33').Value = 'Especial'
$ws.Range("M33").Value = 80
$ws.Range("N33").Value = 12000
$ws.Range("O33").Value = 12000
$ws.Range("P33").Value = 12000
$ws.Range('Q33').Value = '$/caja 18 kilos empedrada'
$ws.Range('R33').Value = 'Región del Maule'
$ws.Range("S33").Value = 667
$ws.Range("T33").Value = 18

# Row 34
$ws.Range("A34").Value = 7
$ws.Range('B34').Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range('C34').Value = 'Ñuble'
$ws.Range("D34").Value = (Get-Date -Year 2023 -Month 6 -Day 16 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E34").Value = 16
$ws.Range('F34').Value = 'Fruta'
$ws.Range("G34").Value = 100104
$ws.Range('H34').Value = 'Frutos de pepita'
$ws.Range("I34").Value = 100104003
$ws.Range('J34').Value = 'Membrillo'
$ws.Range('K34').Value = 'Champion'
$ws.Range('L34').Value = 'Primera'
$ws.Range("M34").Value = 80
$ws.Range("N34").Value = 10000
$ws.Range("O34").Value = 10000
$ws.Range("P34").Value = 10000
$ws.Range('Q34').Value = '$/caja 18 kilos empedrada'
$ws.Range('R34').Value = 'Región del Maule'
$ws.Range("S34").Value = 556
$ws.Range("T34").Value = 18

# Row 35
$ws.Range("A35").Value = 7
$ws.Range('B35').Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range('C35').Value = 'Ñuble'
$ws.Range("D35").Value = (Get-Date -Year 2023 -Month 6 -Day 16 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E35").Value = 16
$ws.Range('F35').Value = 'Fruta'
$ws.Range("G35").Value = 100104
$ws.Range('H35').Value = 'Frutos de pepita'
$ws.Range("I35").Value = 100104003
$ws.Range('J35').Value = 'Membrillo'
$ws.Range('K35').Value = 'Champion'
$ws.Range('L35').Value = 'Segunda'
$ws.Range("M35").Value = 70
$ws.Range("N35").Value = 8000
$ws.Range("O35").Value = 8000
$ws.Range("P35").Value = 8000
$ws.Range('Q35').Value = '$/caja 18 kilos empedrada'
$ws.Range('R35').Value = 'Región del Maule'
$ws.Range("S35").Value = 444
$ws.Range("T35").Value = 18

# Row 36
$ws.Range("A36").Value = 7
$ws.Range('B36').Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range('C36').Value = 'Ñuble'
$ws.Range("D36").Value = (Get-Date -Year 2023 -Month 4 -Day 27 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E36").Value = 16
$ws.Range('F36').Value = 'Fruta'
$ws.Range("G36").Value = 100104
$ws.Range('H36').Value = 'Frutos de pepita'
$ws.Range("I36").Value = 100104003
$ws.Range('J36').Value = 'Membrillo'
$ws.Range('K36').Value = 'Champion'
$ws.Range('L36').Value = 'Especial'
$ws.Range("M36").Value = 40
$ws.Range("N36").Value = 13000
$ws.Range("O36").Value = 13000
$ws.Range("P36").Value = 13000
$ws.Range('Q36').Value = '$/caja 18 kilos empedrada'
$ws.Range('R36').Value = 'Región de O''Higgins'
$ws.Range("S36").Value = 722
$ws.Range("T36").Value = 18

# Row 37
$ws.Range("A37").Value = 7
$ws.Range('B37').Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range('C37').Value = 'Ñuble'
$ws.Range("D37").Value = (Get-Date -Year 2023 -Month 4 -Day 27 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E37").Value = 16
$ws.Range('F37').Value = 'Fruta'
$ws.Range("G37").Value = 100104
$ws.Range('H37').Value = 'Frutos de pepita'
$ws.Range("I37").Value = 100104003
$ws.Range('J37').Value = 'Membrillo'
$ws.Range('K37').Value = 'Champion'
$ws.Range('L37').Value = 'Primera'
$ws.Range("M37").Value = 50
$ws.Range("N37").Value = 12000
$ws.Range("O37").Value = 12000
$ws.Range("P37").Value = 12000
$ws.Range('Q37').Value = '$/caja 18 kilos empedrada'
$ws.Range('R37').Value = 'Región de O''Higgins'
$ws.Range("S37").Value = 667
$ws.Range("T37").Value = 18

# Row 38
$ws.Range("A38").Value = 7
$ws.Range('B38').Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range('C38').Value = 'Ñuble'
$ws.Range("D38").Value = (Get-Date -Year 2023 -Month 6 -Day 13 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E38").Value = 16
$ws.Range('F38').Value = 'Fruta'
$ws.Range("G38").Value = 100104
$ws.Range('H38').Value = 'Frutos de pepita'
$ws.Range("I38").Value = 100104003
$ws.Range('J38').Value = 'Membrillo'
$ws.Range('K38').Value = 'Champion'
$ws.Range('L38').Value = 'Especial'
$ws.Range("M38").Value = 80
$ws.Range("N38").Value = 11000
$ws.Range("O38").Value = 11000
$ws.Range("P38").Value = 11000
$ws.Range('Q38').Value = '$/caja 18 kilos empedrada'
$ws.Range('R38').Value = 'Región del Maule'
$ws.Range("S38").Value = 611
$ws.Range("T38").Value = 18

# Row 39
$ws.Range("A39").Value = 7
$ws.Range('B39').Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range('C39').Value = 'Ñuble'
$ws.Range("D39").Value = (Get-Date -Year 2023 -Month 6 -Day 13 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E39").Value = 16
$ws.Range('F39').Value = 'Fruta'
$ws.Range("G39").Value = 100104
$ws.Range('H39').Value = 'Frutos de pepita'
$ws.Range("I39").Value = 100104003
$ws.Range('J39').Value = 'Membrillo'
$ws.Range('K39').Value = 'Champion'
$ws.Range('L39').Value = 'Primera'
$ws.Range("M39").Value = 50
$ws.Range("N39").Value = 10000
$ws.Range("O39").Value = 10000
$ws.Range("P39").Value = 10000
$ws.Range('Q39').Value = '$/caja 18 kilos empedrada'
$ws.Range('R39').Value = 'Región del Maule'
$ws.Range("S39").Value = 556
$ws.Range("T39").Value = 18

# Row 40
$ws.Range("A40").Value = 7
$ws.Range('B40').Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range('C40').Value = 'Ñuble'
$ws.Range("D40").Value = (Get-Date -Year 2023 -Month 6 -Day 13 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E40").Value = 16
$ws.Range('F40').Value = 'Fruta'
$ws.Range("G40").Value = 100104
$ws.Range('H40').Value = 'Frutos de pepita'
$ws.Range("I40").Value = 100104003
$ws.Range('J40').Value = 'Membrillo'
$ws.Range('K40').Value = 'Champion'
$ws.Range('L40').Value = 'Segunda'
$ws.Range("M40").Value = 40
$ws.Range("N40").Value = 8000
$ws.Range("O40").Value = 8000
$ws.Range("P40").Value = 8000
$ws.Range('Q40').Value = '$/caja 18 kilos empedrada'
$ws.Range('R40').Value = 'Región del Maule'
$ws.Range("S40").Value = 444
$ws.Range("T40").Value = 18

# Row 41
$ws.Range("A41").Value = 7
$ws.Range('B41').Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range('C41').Value = 'Ñuble'
$ws.Range("D41").NumberFormat = $ws.Range("D39").NumberFormat
$ws.Range("D41").Value = (Get-Date -Year 2023 -Month 6 -Day 8 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E41").Value = 16
$ws.Range('F41').Value = 'Fruta'
$ws.Range("G41").Value = 100104
$ws.Range('H41').Value = 'Frutos de pepita'
$ws.Range("I41").Value = 100104003
$ws.Range('J41').Value = 'Membrillo'
$ws.Range('K41').Value = 'Champion'
$ws.Range('L41').Value = 'Primera'
$ws.Range("M41").Value = 50
$ws.Range("N41").Value = 10000
$ws.Range("O41").Value = 10000
$ws.Range("P41").Value = 10000
$ws.Range('Q41').Value = '$/caja 18 kilos empedrada'
$ws.Range('R41').Value = 'Región del Maule'
$ws.Range("S41").Value = 556
$ws.Range("T41").Value = 18

# Row 42
$ws.Range("A42").Value = 7
$ws.Range('B42').Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range('C42').Value = 'Ñuble'
$ws.Range("D42").NumberFormat = $ws.Range("D39").NumberFormat
$ws.Range("D42").Value = (Get-Date -Year 2023 -Month 4 -Day 5 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E42").Value = 16
$ws.Range('F42').Value = 'Fruta'
$ws.Range("G42").Value = 100104
$ws.Range('H42').Value = 'Frutos de pepita'
$ws.Range("I42").Value = 100104003
$ws.Range('J42').Value = 'Membrillo'
$ws.Range('K42').Value = 'Champion'
$ws.Range('L42').Value = 'Primera'
$ws.Range("M42").Value = 50
$ws.Range("N42").Value = 12000
$ws.Range("O42").Value = 12000
$ws.Range("P42").Value = 12000
$ws.Range('Q42').Value = '$/caja 18 kilos granel'
$ws.Range('R42').Value = 'Región de O''Higgins'
$ws.Range("S42").Value = 667
$ws.Range("T42").Value = 18

# Row 43
$ws.Range("A43").Value = 7
$ws.Range('B43').Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range('C43').Value = 'Ñuble'
$ws.Range("D43").NumberFormat = $ws.Range("D39").NumberFormat
$ws.Range("D43").Value = (Get-Date -Year 2023 -Month 6 -Day 6 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E43").Value = 16
$ws.Range('F43').Value = 'Fruta'
$ws.Range("G43").Value = 100104
$ws.Range('H43').Value = 'Frutos de pepita'
$ws.Range("I43").Value = 100104003
$ws.Range('J43').Value = 'Membrillo'
$ws.Range('K43').Value = 'Champion'
$ws.Range('L43').Value = 'Primera'
$ws.Range("M43").Value = 60
$ws.Range("N43").Value = 9000
$ws.Range("O43").Value = 10000
$ws.Range("P43").Value = 9500
$ws.Range('Q43').Value = '$/caja 18 kilos empedrada'
$ws.Range('R43').Value = 'Región del Maule'
$ws.Range("S43").Value = 528
$ws.Range("T43").Value = 18
